# Updated cryptos list with latest price/volume data (GitHub Actions scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.203.57"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.815.93"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("D15").Value = "3.255.35"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "2.801.27"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.902"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "52.190.74"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  +6.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").Value = "0.0₃0984"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0453"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +31.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0837"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  +9.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "2.062.86"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.945"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.46%  "
